{"js": "// Update the division problems in the practice-sheet table.\n// Each cell contains a unique \"NNN\u00f7N=\" expression; replace the old\n// expression text with the new one, cell by cell, using Body.search.\nconst replacements = [\n  [\"120\u00f72=\", \"464\u00f74=\"],\n  [\"828\u00f76=\", \"404\u00f75=\"],\n  [\"702\u00f75=\", \"150\u00f77=\"],\n  [\"984\u00f78=\", \"540\u00f73=\"],\n  [\"951\u00f72=\", \"445\u00f75=\"],\n  [\"391\u00f79=\", \"101\u00f74=\"],\n  [\"415\u00f72=\", \"518\u00f73=\"],\n  [\"437\u00f79=\", \"414\u00f76=\"],\n  [\"707\u00f78=\", \"303\u00f77=\"],\n  [\"340\u00f73=\", \"402\u00f75=\"],\n  [\"808\u00f73=\", \"170\u00f76=\"],\n  [\"511\u00f78=\", \"946\u00f78=\"],\n  [\"248\u00f79=\", \"527\u00f73=\"],\n  [\"734\u00f73=\", \"291\u00f73=\"],\n  [\"502\u00f75=\", \"184\u00f79=\"],\n  [\"857\u00f73=\", \"121\u00f78=\"],\n  [\"495\u00f78=\", \"201\u00f78=\"],\n  [\"179\u00f74=\", \"844\u00f74=\"],\n  [\"870\u00f77=\", \"812\u00f74=\"],\n  [\"749\u00f73=\", \"236\u00f74=\"],\n  [\"921\u00f73=\", \"156\u00f75=\"],\n  [\"855\u00f73=\", \"905\u00f79=\"],\n  [\"772\u00f72=\", \"392\u00f76=\"],\n  [\"376\u00f76=\", \"743\u00f76=\"],\n  [\"506\u00f76=\", \"629\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the practice-sheet table.\n# Each cell contains a unique \"NNN\u00f7N=\" expression; find/replace the old\n# expression text with the new one throughout the document.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"120\u00f72=\", \"464\u00f74=\"),\n    @(\"828\u00f76=\", \"404\u00f75=\"),\n    @(\"702\u00f75=\", \"150\u00f77=\"),\n    @(\"984\u00f78=\", \"540\u00f73=\"),\n    @(\"951\u00f72=\", \"445\u00f75=\"),\n    @(\"391\u00f79=\", \"101\u00f74=\"),\n    @(\"415\u00f72=\", \"518\u00f73=\"),\n    @(\"437\u00f79=\", \"414\u00f76=\"),\n    @(\"707\u00f78=\", \"303\u00f77=\"),\n    @(\"340\u00f73=\", \"402\u00f75=\"),\n    @(\"808\u00f73=\", \"170\u00f76=\"),\n    @(\"511\u00f78=\", \"946\u00f78=\"),\n    @(\"248\u00f79=\", \"527\u00f73=\"),\n    @(\"734\u00f73=\", \"291\u00f73=\"),\n    @(\"502\u00f75=\", \"184\u00f79=\"),\n    @(\"857\u00f73=\", \"121\u00f78=\"),\n    @(\"495\u00f78=\", \"201\u00f78=\"),\n    @(\"179\u00f74=\", \"844\u00f74=\"),\n    @(\"870\u00f77=\", \"812\u00f74=\"),\n    @(\"749\u00f73=\", \"236\u00f74=\"),\n    @(\"921\u00f73=\", \"156\u00f75=\"),\n    @(\"855\u00f73=\", \"905\u00f79=\"),\n    @(\"772\u00f72=\", \"392\u00f76=\"),\n    @(\"376\u00f76=\", \"743\u00f76=\"),\n    @(\"506\u00f76=\", \"629\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
